$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: populate the brand-new rows 8 and 9 (minus the "Datasheet" column,
#     which is filled in later alongside row 7's edits) ---
$ws.Range("A8").Value2 = "Reflective Optical Sensor with Transistor Output"
$ws.Range("C8").Value2 = ".99TL"
$ws.Range("D8").Value2 = "kalsörde"
$ws.Range("E8").Value2 = "https://www.robotistan.com/tcrt5000-kizilotesi-sensor?query=k%C4%B1z%C4%B1l%C3%B6tesi%20sens%C3%B6r&"
$ws.Range("F8").Value2 = "IR sensör for 0.2 to 15mm"

$ws.Range("A9").Value2 = "IR sensor card (8)"
$ws.Range("C9").Value2 = "67TL"
$ws.Range("E9").Value2 = "https://www.robotistan.com/qtr-8rc-kizilotesi-sensor?query=k%C4%B1z%C4%B1l%C3%B6tesi%20sens%C3%B6r&"
$ws.Range("F9").Value2 = "IR sensör for 3mm"

# --- Step 2: fill the "Datasheet" column for rows 7 and 9 with a placeholder,
#     then fix the component name typo in row 7 ("Lase" -> "Laser") ---
$ws.Range("D7").Value2 = "_"
$ws.Range("D9").Value2 = "_"
$ws.Range("A7").Value2 = "waveshare Laser sensor"

# --- Step 3: add the new component in row 10 ---
$ws.Range("A10").Value2 = "TF mini LiDAR"
$ws.Range("C10").Value2 = "273.83TL"
$ws.Range("D10").Value2 = "klasörde"
$ws.Range("E10").Value2 = "https://www.direnc.net/tf-mini-lidar-tof-lazer-menzil-sensoru-dfrobot"
$ws.Range("F10").Value2 = "Laser menzil sensörü for mapping"

# --- Step 4: add hyperlinks for the Web Link cells of the new rows, matching
#     the style already used by E7 ---
$ws.Hyperlinks.Add($ws.Range("E8"), $ws.Range("E8").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E9"), $ws.Range("E9").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E10"), $ws.Range("E10").Value2) | Out-Null

$ws.Range("E8").Style = $ws.Range("E7").Style
$ws.Range("E9").Style = $ws.Range("E7").Style
$ws.Range("E10").Style = $ws.Range("E7").Style

# --- Step 5: move the active selection to F10 ---
$ws.Range("F10").Select() | Out-Null
